$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 6.175892353057861
$ws.Range("B1").Value = 4.755585670471191
$ws.Range("C1").Value = 4.897356510162354
$ws.Range("D1").Value = 4.349461555480957
$ws.Range("E1").Value = 3.091338634490967
